$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 106 (weekly price update), pushing the
# existing rows 106-128 down to 107-129.
$ws.Rows.Item(106).Insert()

# The new row starts as a duplicate of the record that used to be at row 106
# (now shifted to row 107), with the date, volume and origin updated for
# this week's entry.
for ($col = 1; $col -le 20; $col++) {
    $src = $ws.Cells.Item(107, $col)
    $dst = $ws.Cells.Item(106, $col)
    $dst.Value = $src.Value2
}

$ws.Range("D106").Value = 45007
$ws.Range("M106").Value = 100
$ws.Range("R106").Value = "Provincia de Cachapoal"
